$wb = $excel.ActiveWorkbook

# Add a new row of data to the "managecategorypage" sheet (5th sheet).
$ws = $wb.Worksheets.Item("managecategorypage")

$ws.Range("A4").Value = "Perfumes"
$ws.Range("B4").Value = "Yes"
$ws.Range("C4").Value = "Yes"

# Make this the active sheet / active cell, matching the saved view state
# (tabSelected moves from the login page sheet to this one, activeTab/firstSheet
# move to this sheet in the workbook view).
$ws.Activate()
$ws.Range("C4").Select()
